# TC_5-FDR_E2E-FDR-2938-SYSTEM ADMIN REFERENCE DATA
# Update the Global Variables sheet's Draw Date / Pay Date to 10/12/2020 (serial 44116),
# which cascades into all cached formula results on the Voucher sheets that reference them.

$wb = $excel.ActiveWorkbook

$gv = $wb.Worksheets.Item("Global Variables")
$gv.Range("B1").Value = 44116
$gv.Range("B2").Value = 44116

# On Voucher1, change selection away from the full-column selection to B33,
# and it should no longer be the tab that is active/selected when saved.
$v1 = $wb.Worksheets.Item("Voucher1")
$v1.Activate()
$v1.Range("B33").Select()

# Make "Global Variables" the active/selected sheet, keeping its existing
# selection at B3.
$gv.Activate()
$gv.Range("B3").Select()

$wb.Save()
